# Applies the marksheet auto-grading update:
#  - fills in the summary rows (10-12) with the computed right/wrong/not-attempted/max
#    counts and the resulting total/percentage
#  - removes the third "Student Ans / Correct Ans" question block (columns G:H)
#  - trims the second question block (columns D:E) down to the 3 questions that
#    were actually asked
#  - fills the "Student Ans" columns (A and D) with the grading result: the
#    correct answer text when the student's answer matched (shown with the
#    correctStyle look), the student's (wrong) answer when it did not match
#    (incorrectStyle), or leaves the cell blank when the question was not
#    attempted (normalStyle, unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary block (rows 10-12) ----------------------------------------
# Row 10/11/12 labels pick up the same "mtitleStyle" look already used by A9.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("B10").Value = 26
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

$ws.Range("B12").Value = 104
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "103/112"

# ---- Drop the third question block (columns G:H) ------------------------
$ws.Range("G15:H40").Clear()

# ---- Trim the second question block (columns D:E) to 3 questions --------
$ws.Range("D19:E40").Clear()

# ---- Grade the first question block (column A vs column B) --------------
# For every row, compare the "Correct Ans" (column B) against what should be
# shown as the "Student Ans" (column A). A match is shown in correctStyle
# (sampled from B10), a mismatch in incorrectStyle (sampled from C10); rows
# left blank were not attempted and keep their original normalStyle look.
$studentAnsA = @{
    16 = "Option A"
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    20 = "Option B"
    21 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = "Option A"
    25 = "Option A"
    27 = "Option A"
    28 = "Option D"
    29 = "Option D"
    30 = "Option B"
    31 = "Option C"
    32 = "Option C"
    33 = "Option D"
    34 = "Option B"
    35 = "Option D"
    36 = "Option A"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
    40 = "Option D"
}

foreach ($row in 16..40) {
    $cellA = $ws.Cells.Item($row, 1)
    if (-not $studentAnsA.ContainsKey($row)) {
        # Row 26: left as not attempted - no change needed.
        continue
    }
    $answer = $studentAnsA[$row]
    $correct = $ws.Cells.Item($row, 2).Value
    if ($answer -eq $correct) {
        $ws.Range("B10").Copy()
    } else {
        $ws.Range("C10").Copy()
    }
    $cellA.PasteSpecial(-4122)
    $cellA.Value = $answer
}

# ---- Grade the second question block (column D vs column E) -------------
foreach ($row in 16..18) {
    $cellD = $ws.Cells.Item($row, 4)
    $correct = $ws.Cells.Item($row, 5).Value
    $ws.Range("B10").Copy()
    $cellD.PasteSpecial(-4122)
    $cellD.Value = $correct
}
